$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: Cells.Clear() also drops now-unused shared strings (e.g. the old
# generic "Name"/"Description" header and several stale descriptions) so the
# saved workbook does not retain orphaned sharedStrings entries.
$ws.Cells.Clear()

$data = @(
    @('Feature Name', 'Feature Description'),
    @('correct', 'Whether or not the question was answered correctly'),
    @('questionID', 'Unique question identifier'),
    @('examID', 'Unique exam identifier'),
    @('course', 'Either P, FM, or MFE'),
    @('exam_type', 'Either quiz ("q") or exam ("e")'),
    @('creation_dt', 'Date on which the exam or quiz was created'),
    @('creation_hr', 'Hour of creation'),
    @('marked', 'Adapt "mark question" option'),
    @('q_ordinal', 'The order in which the question appeared in the exam (e.g., from 1-30 for exams)'),
    @('difficulty', 'The Adapt-generated difficulty from 0-11'),
    @('nth_exam', 'Exam count number'),
    @('nth_e_or_q', 'Exam or quiz count number'),
    @('weekday', 'Day of week'),
    @('minutes_used', 'Minutes of screen view time'),
    @('cat1', 'Major section category, limited to the most frequent 10.  All else are "other"'),
    @('cat2', 'Major section category, limited to the most frequent 10.  All else are "other"'),
    @('cat3', 'Major section category, limited to the most frequent 10.  All else are "other"'),
    @('subcat1', 'Minor section category, limited to the most frequent 10.  All else are "other"'),
    @('subcat2', 'Minor section category, limited to the most frequent 10.  All else are "other"'),
    @('subcat3', 'Minor section category, limited to the most frequent 10.  All else are "other"'),
    @('approx_time_remaining', 'Approximate remaining time on the exam.  This is not "real" because I would almost always skip around from question to question'),
    @('EL_begin', 'Earned level at start of exam'),
    @('EL_change', 'Change in Earned level from start-finish of exam'),
    @('EL_end', 'Ending Earned Level'),
    @('hrs_since_previous_e', 'Hours since the most recent exam'),
    @('hist_subcat_n', 'The number of historical questions in the current question category'),
    @('hist_subcat_diff', 'The sum of difficuly for historical question in the current question subcategory'),
    @('hist_subcat_diff_correct', 'Subset of hist_subcat_diff which were answered correctly'),
    @('hist_total_time_e', 'The historical sum of minutes_used for correct questions in the current subcategory for exams'),
    @('hist_total_time_e_correct', 'The subset of hist_total_time_e which were correct'),
    @('hist_total_time_q', 'Historical total quiz time for questions in current subcategory'),
    @('hist_total_time_e', 'hist_total_time_q but for exams'),
    @('hist_greater_diff_time', 'Total time spent on exams in current subcategory which have greater difficulty'),
    @('hist_greater_diff_time_correct', 'hist_greater_diff_time which were correct'),
    @('hist_greater_diff', 'same as hist_greater_diff_time but summing over difficulty'),
    @('hist_net_diff', 'same as hist_greater_diff but adds for correct questions and subtracts difficulty for incorrect'),
    @('hist_n_marked', 'Number of marked questions for exams in the current subcategory'),
    @('hist_n_marked_incorrect', 'Same as hist_n_marked but only for incorrect questions'),
    @('hist_repeat_question', 'The total number of times that the current question has already been seen in exams or quizzes'),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Re-apply header formatting (bold row 1) that Clear() wiped out
$ws.Range("A1:B1").Font.Bold = $true

# Widen column A for the longer "Feature Name" header / values
$ws.Columns.Item(1).ColumnWidth = 30.8

# Restore view state: scrolled so row 13 is at the top, B18 selected
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$ws.Range("B18").Select()

